$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the row above into the new row 11
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the values for the new row
$ws.Range("A11").Value = "biginteger"
$ws.Range("B11").Value = 125702689
